# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.227.96"
$ws.Range("E2").Value = "  -0.59%  "

# Row 3
$ws.Range("D3").Value = "2.238.91"
$ws.Range("E3").Value = "  -0.09%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.88"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.37%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.14"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.07%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  -3.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.94%  "

# Row 11
$ws.Range("E11").Value = "  -0.13%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.92%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.28%  "

# Row 14
$ws.Range("D14").Value = "2.568.57"
$ws.Range("E14").Value = "  -0.04%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.839"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.95%  "

# Row 17
$ws.Range("D17").Value = "2.194.17"
$ws.Range("E17").Value = "  -2.23%  "

# Row 18
$ws.Range("D18").Value = "42.099.42"
$ws.Range("E18").Value = "  -0.51%  "

# Row 19
$ws.Range("E19").Value = "  -4.23%  "

# Row 20
$ws.Range("E20").Value = "  +0.64%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.11%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.63%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.57%  "

# Row 27
$ws.Range("E27").Value = "  -0.49%  "

# Row 28
$ws.Range("E28").Value = "  -1.14%  "

# Row 29
$ws.Range("E29").Value = "  -2.32%  "

# Row 30
$ws.Range("E30").Value = "  +0.25%  "

# Row 31
$ws.Range("E31").Value = "  -1.52%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.50%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0806"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.95%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.32%  "

# Row 35
$ws.Range("E35").Value = "  -0.50%  "

# Row 36
$ws.Range("E36").Value = "  -6.82%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.89%  "

# Row 38
$ws.Range("E38").Value = "  -2.25%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.92%  "

# Row 40
$ws.Range("E40").Value = "  -1.82%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.63%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.98%  "

# Row 43
$ws.Range("E43").Value = "  -0.79%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.35%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.38%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.100"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.26%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.65%  "

# Row 48
$ws.Range("E48").Value = "  -0.17%  "

# Row 49
$ws.Range("E49").Value = "  -2.76%  "

# Row 50
$ws.Range("E50").Value = "  -1.83%  "

# Row 51
$ws.Range("D51").Value = "2.443.25"
$ws.Range("E51").Value = "  -0.07%  "
